$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hbegf"
$ws.Range("C2").Value = "Erbb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.572976999999999
$ws.Range("H2").Value = 22.718931
$ws.Range("I2").Value = 0.4497670593913077
$ws.Range("J2").Value = 0.4497670593913078
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.179771666666667
$ws.Range("N2").Value = 6.539315
$ws.Range("O2").Value = 0.2349306639444428
$ws.Range("P2").Value = 0.2349306639444428
$ws.Range("Q2").Value = 16.50736069691833
$ws.Range("R2").Value = 148.566246272265
$ws.Range("S2").Value = 0.1056640738831396
$ws.Range("T2").Value = 0.1056640738831396

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hbegf"
$ws.Range("C3").Value = "Erbb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.572976999999999
$ws.Range("H3").Value = 22.718931
$ws.Range("I3").Value = 0.4497670593913077
$ws.Range("J3").Value = 0.4497670593913078
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.237801
$ws.Range("N3").Value = 12.713403
$ws.Range("O3").Value = 0.4567402255103586
$ws.Range("P3").Value = 0.4567402255103586
$ws.Range("Q3").Value = 32.092769503577
$ws.Range("R3").Value = 288.834925532193
$ws.Range("S3").Value = 0.2054267081335167
$ws.Range("T3").Value = 0.2054267081335168

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hbegf"
$ws.Range("C4").Value = "Erbb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.572976999999999
$ws.Range("H4").Value = 22.718931
$ws.Range("I4").Value = 0.4497670593913077
$ws.Range("J4").Value = 0.4497670593913078
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.860789
$ws.Range("N4").Value = 8.582367000000001
$ws.Range("O4").Value = 0.3083291105451987
$ws.Range("P4").Value = 0.3083291105451986
$ws.Range("Q4").Value = 21.664689298853
$ws.Range("R4").Value = 194.982203689677
$ws.Range("S4").Value = 0.1386762773746515
$ws.Range("T4").Value = 0.1386762773746515

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hbegf"
$ws.Range("C5").Value = "Erbb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.784025666666667
$ws.Range("H5").Value = 8.352077
$ws.Range("I5").Value = 0.1653462089435359
$ws.Range("J5").Value = 0.1653462089435359
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.179771666666667
$ws.Range("N5").Value = 6.539315
$ws.Range("O5").Value = 0.2349306639444428
$ws.Range("P5").Value = 0.2349306639444428
$ws.Range("Q5").Value = 6.068540267472778
$ws.Range("R5").Value = 54.616862407255
$ws.Range("S5").Value = 0.03884489464780146
$ws.Range("T5").Value = 0.03884489464780146

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hbegf"
$ws.Range("C6").Value = "Erbb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.784025666666667
$ws.Range("H6").Value = 8.352077
$ws.Range("I6").Value = 0.1653462089435359
$ws.Range("J6").Value = 0.1653462089435359
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.237801
$ws.Range("N6").Value = 12.713403
$ws.Range("O6").Value = 0.4567402255103586
$ws.Range("P6").Value = 0.4567402255103586
$ws.Range("Q6").Value = 11.79814675422567
$ws.Range("R6").Value = 106.183320788031
$ws.Range("S6").Value = 0.07552026476015347
$ws.Range("T6").Value = 0.07552026476015347

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hbegf"
$ws.Range("C7").Value = "Erbb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.784025666666667
$ws.Range("H7").Value = 8.352077
$ws.Range("I7").Value = 0.1653462089435359
$ws.Range("J7").Value = 0.1653462089435359
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.860789
$ws.Range("N7").Value = 8.582367000000001
$ws.Range("O7").Value = 0.3083291105451987
$ws.Range("P7").Value = 0.3083291105451986
$ws.Range("Q7").Value = 7.964510002917668
$ws.Range("R7").Value = 71.68059002625901
$ws.Range("S7").Value = 0.05098104953558101
$ws.Range("T7").Value = 0.050981049535581

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Hbegf"
$ws.Range("C8").Value = "Erbb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6.480551000000001
$ws.Range("H8").Value = 19.441653
$ws.Range("I8").Value = 0.3848867316651562
$ws.Range("J8").Value = 0.3848867316651562
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.179771666666667
$ws.Range("N8").Value = 6.539315
$ws.Range("O8").Value = 0.2349306639444428
$ws.Range("P8").Value = 0.2349306639444428
$ws.Range("Q8").Value = 14.12612145418834
$ws.Range("R8").Value = 127.135093087695
$ws.Range("S8").Value = 0.09042169541350173
$ws.Range("T8").Value = 0.09042169541350174

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Hbegf"
$ws.Range("C9").Value = "Erbb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6.480551000000001
$ws.Range("H9").Value = 19.441653
$ws.Range("I9").Value = 0.3848867316651562
$ws.Range("J9").Value = 0.3848867316651562
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.237801
$ws.Range("N9").Value = 12.713403
$ws.Range("O9").Value = 0.4567402255103586
$ws.Range("P9").Value = 0.4567402255103586
$ws.Range("Q9").Value = 27.46328550835101
$ws.Range("R9").Value = 247.169569575159
$ws.Range("S9").Value = 0.1757932526166883
$ws.Range("T9").Value = 0.1757932526166884

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Hbegf"
$ws.Range("C10").Value = "Erbb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 6.480551000000001
$ws.Range("H10").Value = 19.441653
$ws.Range("I10").Value = 0.3848867316651562
$ws.Range("J10").Value = 0.3848867316651562
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.860789
$ws.Range("N10").Value = 8.582367000000001
$ws.Range("O10").Value = 0.3083291105451987
$ws.Range("P10").Value = 0.3083291105451986
$ws.Range("Q10").Value = 18.539489014739
$ws.Range("R10").Value = 166.8554011326511
$ws.Range("S10").Value = 0.1186717836349662
$ws.Range("T10").Value = 0.1186717836349662

